$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos data (prices and volume %) scraped on 2024-03-03.
# Force text format on target cells so numeric-looking strings (e.g. "0.999",
# "62.787.77") are stored as text, matching the original inline-string layout.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.787.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.468.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.61"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.65"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.023.83"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.140"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.448.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.703.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "462.38"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.70"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +17.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.94%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "147.59"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.321"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.69"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.34"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +13.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₃0560"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +30.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.16"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.140"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.96%  "
